# Apply passive-trial value tweaks ("Hjemme passive tweaks lichtwark deleted values")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (trial/length headers)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2
$ws.Range("B2").Value = 45.976771617899388
$ws.Range("C2").Value = 53.906554568370282
$ws.Range("D2").Value = 48.79964776451186
$ws.Range("E2").Value = 54.878858733544178

# Row 3
$ws.Range("B3").Value = 42.602687163599157
$ws.Range("C3").Value = 44.255729989544079
$ws.Range("D3").Value = 43.767603607131896
$ws.Range("E3").Value = 54.559152646880094

# Update the active selection to match the edited range
$ws.Range("B1:E3").Select()
